$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.095.85'
$ws.Range("E2").Value = '  -3.02%  '

$ws.Range("D3").Value = '1.866.27'
$ws.Range("E3").Value = '  -2.31%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = '''307.05'
$ws.Range("E5").Value = '  -1.94%  '

$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").Value = '''0.5090'
$ws.Range("E7").Value = '  +1.59%  '

$ws.Range("D8").Value = '''0.3735'
$ws.Range("E8").Value = '  -2.23%  '

$ws.Range("D9").Value = '''0.07141'
$ws.Range("E9").Value = '  -2.48%  '

$ws.Range("D10").Value = '''0.8869'
$ws.Range("E10").Value = '  -2.72%  '

$ws.Range("D11").Value = '''20.57'
$ws.Range("E11").Value = '  -3.25%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '''0.07538'
$ws.Range("E12").Value = '  -1.76%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.855.30'
$ws.Range("E13").Value = '  -2.97%  '

$ws.Range("D14").Value = '''5.314'
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").Value = '''89.28'
$ws.Range("E15").Value = '  -3.96%  '

$ws.Range("D16").Value = '''1.001'

$ws.Range("D17").Value = '''0.000008455'
$ws.Range("E17").Value = '  -3.34%  '

$ws.Range("D18").Value = '''14.10'
$ws.Range("E18").Value = '  -4.16%  '

$ws.Range("D19").Value = '''0.9999'
$ws.Range("E19").Value = '  +0.09%  '

$ws.Range("D20").Value = '27.126.90'
$ws.Range("E20").Value = '  -3.02%  '

$ws.Range("D21").Value = '''5.061'
$ws.Range("E21").Value = '  -2.38%  '

$ws.Range("D22").Value = '2.111.01'
$ws.Range("E22").Value = '  -2.73%  '

$ws.Range("D23").Value = '''10.55'
$ws.Range("E23").Value = '  -2.89%  '

$ws.Range("D24").Value = '''6.483'
$ws.Range("E24").Value = '  -1.98%  '

$ws.Range("D25").Value = '''150.23'
$ws.Range("E25").Value = '  -1.92%  '

$ws.Range("D26").Value = '''1.843'
$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  -2.79%  '

$ws.Range("D28").Value = '''2.099'
$ws.Range("E28").Value = '  -5.16%  '

$ws.Range("D29").Value = '''112.60'
$ws.Range("E29").Value = '  -2.50%  '

$ws.Range("D30").Value = '''4.749'
$ws.Range("E30").Value = '  -3.73%  '

$ws.Range("D31").Value = '''4.680'
$ws.Range("E31").Value = '  -3.83%  '

$ws.Range("D32").Value = '''0.09046'
$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("E33").Value = '  -3.09%  '

$ws.Range("D34").Value = '''3.093'
$ws.Range("E34").Value = '  -3.64%  '

$ws.Range("D35").Value = '''1.160'
$ws.Range("E35").Value = '  -6.21%  '

$ws.Range("D36").Value = '''0.7352'
$ws.Range("E36").Value = '  -5.55%  '

$ws.Range("D37").Value = '''0.02038'
$ws.Range("E37").Value = '  -2.37%  '

$ws.Range("D38").Value = '''2.484'
$ws.Range("E38").Value = '  -4.45%  '

$ws.Range("D39").Value = '''3.041'
$ws.Range("E39").Value = '  -0.90%  '

$ws.Range("D40").Value = '''1.077'
$ws.Range("E40").Value = '  -1.57%  '

$ws.Range("D41").Value = '''0.5324'
$ws.Range("E41").Value = '  -4.19%  '

$ws.Range("D42").Value = '''6.602'
$ws.Range("E42").Value = '  -4.18%  '

$ws.Range("D43").Value = '''115.95'
$ws.Range("E43").Value = '  +2.04%  '

$ws.Range("D44").Value = '''8.335'
$ws.Range("E44").Value = '  -2.32%  '

$ws.Range("E45").Value = '  -3.11%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''0.9994'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.4625'
$ws.Range("E47").Value = '  -4.49%  '

$ws.Range("D48").Value = '''9.965'
$ws.Range("E48").Value = '  -6.30%  '

$ws.Range("E49").Value = '  -4.76%  '

$ws.Range("D50").Value = '''64.50'
$ws.Range("E50").Value = '  -4.74%  '

$ws.Range("E51").Value = '  -2.14%  '
